$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 129 (shifts existing rows 129..251 down to 130..252)
$ws.Rows.Item(129).Insert()

# Populate the new row 129 with a copy of the constant columns plus the new data
$ws.Range("A129").Value = 5
$ws.Range("B129").Value = "Macroferia Regional de Talca"
$ws.Range("C129").Value = "Maule"
$ws.Range("D129").Value = 44589
$ws.Range("D129").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E129").Value = 7
$ws.Range("F129").Value = 100114014
$ws.Range("G129").Value = "Betarraga"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 4000
$ws.Range("K129").Value = 700
$ws.Range("L129").Value = 700
$ws.Range("M129").Value = 700
$ws.Range("N129").Value = "`$/paquete 5 unidades"
$ws.Range("O129").Value = "Región del Maule"
$ws.Range("P129").Value = 140
$ws.Range("Q129").Value = 5
$ws.Range("R129").Value = "Hortaliza"
